# Add a "Github username" column (column C) to the group info sheet,
# populating it with each member's Github username.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> Github username (row 1 is the header row).
$usernames = [ordered]@{
    1  = "Github username"
    2  = "shailyacharya"
    3  = "sbrown5x"
    4  = "neeldesai01"
    5  = "yc577"
    6  = "yashdhuldhoya"
    7  = "felipegermanos"
    8  = "MoatazGU"
    9  = "MarineAntonio"
    10 = "FegorEO"
    11 = "abigailorbe"
    12 = "geena-panzitta"
    13 = "bp557"
    14 = "AaronShtilerman"
    15 = "MiglePetrou"
    17 = "Benjamin-Tu"
    16 = "f2pHgty8hw"
}

foreach ($row in $usernames.Keys) {
    # Match column C's look (border/fill/font) to column B on the same row
    # before writing the value, same as the rest of the table.
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $usernames[$row]
}

# Give the new column a sensible custom width, same as the other columns.
$ws.Columns.Item(3).ColumnWidth = 17.8

# Reflect the cell the author ended up leaving selected.
[void]$ws.Range("C19").Select()
